$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("MAIN")

# New rows 44-48: filing metadata rows appended to the tracking table.
# Columns: C = Entity name, D = Form type, E = Filing date, H = Filing URL (hyperlink),
# I = Extraction Status, J = Formatting status.

$ws.Cells.Item(44, 3).Value = "Main Street Capital CORP"
$ws.Cells.Item(44, 4).Value = "10-Q"
$ws.Cells.Item(44, 5).Value = 44869
$ws.Cells.Item(44, 5).NumberFormat = "mm-dd-yy"
$ws.Cells.Item(44, 8).Value = "https://www.sec.gov/Archives/edgar/data/1396440/000139644022000034/main-20220930.htm"
$ws.Cells.Item(44, 9).Value = "Done"
$ws.Cells.Item(44, 10).Value = "Need to Check"

$ws.Cells.Item(45, 3).Value = "Main Street Capital CORP"
$ws.Cells.Item(45, 4).Value = "10-K"
$ws.Cells.Item(45, 5).Value = 44981
$ws.Cells.Item(45, 5).NumberFormat = "mm-dd-yy"
$ws.Cells.Item(45, 8).Value = "https://www.sec.gov/Archives/edgar/data/1396440/000139644023000035/main-20221231.htm"
$ws.Cells.Item(45, 9).Value = "Done"
$ws.Cells.Item(45, 10).Value = "Need to Check"

$ws.Cells.Item(46, 3).Value = "Main Street Capital CORP"
$ws.Cells.Item(46, 4).Value = "10-Q"
$ws.Cells.Item(46, 5).Value = 45051
$ws.Cells.Item(46, 5).NumberFormat = "mm-dd-yy"
$ws.Cells.Item(46, 8).Value = "https://www.sec.gov/Archives/edgar/data/1396440/000139644023000067/main-20230331.htm"
$ws.Hyperlinks.Add($ws.Cells.Item(46, 8), "https://www.sec.gov/Archives/edgar/data/1396440/000139644023000067/main-20230331.htm") | Out-Null
$ws.Cells.Item(46, 9).Value = "Done"
$ws.Cells.Item(46, 10).Value = "Need to Check"

$ws.Cells.Item(47, 3).Value = "Main Street Capital CORP"
$ws.Cells.Item(47, 4).Value = "10-Q"
$ws.Cells.Item(47, 5).Value = 45142
$ws.Cells.Item(47, 5).NumberFormat = "mm-dd-yy"
$ws.Cells.Item(47, 8).Value = "https://www.sec.gov/Archives/edgar/data/1396440/000139644023000102/main-20230630.htm"
$ws.Cells.Item(47, 9).Value = "Done"
$ws.Cells.Item(47, 10).Value = "Need to Check"

$ws.Cells.Item(48, 3).Value = "Main Street Capital CORP"
$ws.Cells.Item(48, 4).Value = "10-Q"
$ws.Cells.Item(48, 5).Value = 45233
$ws.Cells.Item(48, 5).NumberFormat = "mm-dd-yy"
$ws.Cells.Item(48, 8).Value = "https://www.sec.gov/Archives/edgar/data/1396440/000139644023000140/main-20230930.htm"
$ws.Cells.Item(48, 9).Value = "Done"
$ws.Cells.Item(48, 10).Value = "Need to Check"

$ws.Range("E51").Select() | Out-Null
